$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing descriptions for Grad_premium and Grad_share rows
# (entered in this order so the shared-strings table matches the source file)
$ws.Range("B41").Value = "The increase in average wages that graduate students can expect having earned a degree from a graduate program"
$ws.Range("B40").Value = "Graduate students as a share of the total"

# Update the view/selection to reflect scrolling down to the newly completed rows
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("B42").Select()
